$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# --- Bad Drivers section ---
$ws.Range("D3").Value = 96.90000000000001
$ws.Range("C4").Value = 357
$ws.Range("D5").Value = 97.8
$ws.Range("C6").Value = 420

# --- Good Drivers section ---
$ws.Range("B16").Value = 56069
$ws.Range("B17").Value = 449371
$ws.Range("B22").Value = 276086
$ws.Range("B23").Value = 625298
$ws.Range("B28").Value = 331283
$ws.Range("B30").Value = 453652
$ws.Range("B38").Value = 96091
$ws.Range("B41").Value = 99549
$ws.Range("B42").Value = 77999
$ws.Range("B46").Value = 175767
$ws.Range("B47").Value = 240182
$ws.Range("B56").Value = 684728
$ws.Range("B58").Value = 210188
$ws.Range("B63").Value = 308481
$ws.Range("B70").Value = 443223
$ws.Range("B72").Value = 109665
$ws.Range("B74").Value = 62515
